# Insert two new data rows (rows 277 and 278) right after existing row 276,
# pushing the existing rows 277..326 down to 279..328.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A277:A278").EntireRow.Insert()

# Row 277: Comercializadora del Agro de Limarí - Nectarín - Artic Sprite - Primera
$ws.Cells.Item(277, 1).Value = 2
$ws.Cells.Item(277, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(277, 3).Value = "Coquimbo"
$ws.Cells.Item(277, 4).Value = 45015
$ws.Cells.Item(277, 5).Value = 4
$ws.Cells.Item(277, 6).Value = "Fruta"
$ws.Cells.Item(277, 7).Value = 100103
$ws.Cells.Item(277, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(277, 9).Value = 100103006
$ws.Cells.Item(277, 10).Value = "Nectarín"
$ws.Cells.Item(277, 11).Value = "Artic Sprite"
$ws.Cells.Item(277, 12).Value = "Primera"
$ws.Cells.Item(277, 13).Value = 14
$ws.Cells.Item(277, 14).Value = 450000
$ws.Cells.Item(277, 15).Value = 460000
$ws.Cells.Item(277, 16).Value = 455000
$ws.Cells.Item(277, 17).Value = "$/bins (420 kilos)"
$ws.Cells.Item(277, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(277, 19).Value = 1083
$ws.Cells.Item(277, 20).Value = 420

# Row 278: Comercializadora del Agro de Limarí - Nectarín - Artic Sprite - Segunda
$ws.Cells.Item(278, 1).Value = 2
$ws.Cells.Item(278, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(278, 3).Value = "Coquimbo"
$ws.Cells.Item(278, 4).Value = 45015
$ws.Cells.Item(278, 5).Value = 4
$ws.Cells.Item(278, 6).Value = "Fruta"
$ws.Cells.Item(278, 7).Value = 100103
$ws.Cells.Item(278, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(278, 9).Value = 100103006
$ws.Cells.Item(278, 10).Value = "Nectarín"
$ws.Cells.Item(278, 11).Value = "Artic Sprite"
$ws.Cells.Item(278, 12).Value = "Segunda"
$ws.Cells.Item(278, 13).Value = 12
$ws.Cells.Item(278, 14).Value = 400000
$ws.Cells.Item(278, 15).Value = 410000
$ws.Cells.Item(278, 16).Value = 405000
$ws.Cells.Item(278, 17).Value = "$/bins (420 kilos)"
$ws.Cells.Item(278, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(278, 19).Value = 964
$ws.Cells.Item(278, 20).Value = 420
